$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the cells we touch so numeric-looking
# strings like "1.002" or "242.75" are preserved verbatim as text
# (matching the source data's inline-string convention) instead of
# being auto-converted to numbers by Excel's input parser.
$targets = @(
  "D2"
  "E2"
  "D3"
  "E3"
  "D4"
  "E4"
  "D5"
  "E5"
  "D6"
  "E6"
  "E7"
  "D8"
  "E8"
  "D9"
  "E9"
  "D10"
  "E10"
  "D11"
  "E11"
  "D12"
  "E12"
  "D13"
  "E13"
  "D14"
  "E14"
  "D15"
  "E15"
  "D16"
  "E16"
  "D17"
  "E17"
  "D18"
  "E18"
  "D19"
  "E19"
  "D20"
  "E20"
  "E21"
  "D22"
  "E22"
  "D23"
  "E23"
  "D24"
  "E24"
  "D25"
  "E25"
  "D26"
  "E26"
  "D27"
  "E27"
  "D28"
  "E28"
  "D29"
  "E29"
  "D30"
  "E30"
  "D31"
  "E31"
  "D32"
  "E32"
  "E33"
  "D34"
  "E34"
  "D35"
  "E35"
  "D36"
  "E36"
  "D37"
  "E37"
  "D38"
  "E38"
  "D39"
  "E39"
  "D40"
  "E40"
  "D41"
  "E41"
  "B42"
  "C42"
  "D42"
  "E42"
  "B43"
  "C43"
  "D43"
  "E43"
  "D44"
  "E44"
  "B45"
  "C45"
  "D45"
  "E45"
  "B46"
  "C46"
  "D46"
  "E46"
  "D47"
  "E47"
  "D48"
  "D49"
  "E49"
  "D50"
  "E50"
  "D51"
  "E51"
)
foreach ($ref in $targets) {
  $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.277.99"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "1.866.37"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "0.7065"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").Value = "242.75"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "0.07845"
$ws.Range("E8").Value = "  -2.90%  "
$ws.Range("D9").Value = "0.3115"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").Value = "24.31"
$ws.Range("E10").Value = "  -4.00%  "
$ws.Range("D11").Value = "0.08028"
$ws.Range("E11").Value = "  -3.85%  "
$ws.Range("D12").Value = "1.901.14"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("D13").Value = "5.196"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("D14").Value = "93.60"
$ws.Range("E14").Value = "  +1.49%  "
$ws.Range("D15").Value = "0.6965"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").Value = "6.349"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "29.648.97"
$ws.Range("E17").Value = "  +0.96%  "
$ws.Range("D18").Value = "0.000008299"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("D19").Value = "252.89"
$ws.Range("E19").Value = "  +4.71%  "
$ws.Range("D20").Value = "2.203.91"
$ws.Range("E20").Value = "  +3.58%  "
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "7.541"
$ws.Range("E23").Value = "  -3.48%  "
$ws.Range("D24").Value = "0.9997"
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "0.1554"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").Value = "9.008"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "160.22"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("D28").Value = "18.74"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("D29").Value = "1.501"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "4.281"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").Value = "4.281"
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("D32").Value = "1.218"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").Value = "1.896"
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").Value = "0.7478"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "1.158"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("D37").Value = "2.708"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "0.01864"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").Value = "1.250.09"
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("D40").Value = "2.749"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "6.262"
$ws.Range("E41").Value = "  -4.64%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "111.63"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.9007"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").Value = "71.99"
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "2.092.17"
$ws.Range("E45").Value = "  +3.43%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "0.00000000126"
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").Value = "1.798"
$ws.Range("D49").Value = "0.5196"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "9.377"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").Value = "1.010"
$ws.Range("E51").Value = "  +1.14%  "
